# "Find Duplicate File in System" whiteboard: jot down the design note and a
# small worked example (index -> path) used while reasoning through the
# approach.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Heading
$ws.Range("F5").Value = "design"

# Worked example: map of index -> path, as scratch notes while tracing
# through a root/a, root/c, root/c/d style directory tree.
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "root/a"

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = "root/c"

$ws.Range("I8").Value = 2
$ws.Range("J8").Value = "root/c/d"

$ws.Range("I9").Value = 3
$ws.Range("J9").Value = "root"

# The path column is the widest entry, so size it to fit its contents.
$ws.Columns.Item(10).AutoFit() | Out-Null

# Leave the cursor where the author's session ended.
$ws.Range("L13").Select() | Out-Null
